# docs: update brew installed software
# Updates the "Mac installed" sheet (Brew column values + a few renamed
# package names) and refreshes the saved selection on "Mac availability"
# and "Mac installed".

$wb = $excel.ActiveWorkbook
$wsInstalled = $wb.Worksheets.Item("Mac installed")
$wsAvailability = $wb.Worksheets.Item("Mac availability")

# --- "Mac installed" sheet: rename a few Package entries to their brew
#     cask/formula ids ---
$wsInstalled.Range("A8").Value  = "visual-studio-code"
$wsInstalled.Range("A9").Value  = "google-chrome"
$wsInstalled.Range("A11").Value = "adobe-acrobat-reader"
$wsInstalled.Range("A18").Value = "android-file-transfer"
$wsInstalled.Range("A19").Value = "epic-games"
$wsInstalled.Range("A24").Value = "logi-options-plus"
$wsInstalled.Range("A26").Value = "karabiner-elements"

# --- "Mac installed" sheet: refresh the "Brew" column (C) now that the
#     installed-via-brew status has been re-checked ---
$wsInstalled.Range("C2").Value  = "✅"
$wsInstalled.Range("C3").Value  = "✅"
$wsInstalled.Range("C4").Value  = "❌"
$wsInstalled.Range("C5").Value  = "✅"
$wsInstalled.Range("C6").Value  = "✅"
$wsInstalled.Range("C7").Value  = "✅"
$wsInstalled.Range("C8").Value  = "✅"
$wsInstalled.Range("C9").Value  = "✅"
$wsInstalled.Range("C10").Value = "✅"
$wsInstalled.Range("C11").Value = "✅"
$wsInstalled.Range("C17").Value = "✅"
$wsInstalled.Range("C18").Value = "✅"
$wsInstalled.Range("C19").Value = "✅"
$wsInstalled.Range("C20").Value = "✅"
$wsInstalled.Range("C21").Value = "❌"
$wsInstalled.Range("C22").Value = "❌"
$wsInstalled.Range("C23").Value = "❌"
$wsInstalled.Range("C24").Value = "(✅)"
$wsInstalled.Range("C25").Value = "✅"
$wsInstalled.Range("C26").Value = "✅"
$wsInstalled.Range("C33").Value = "✅"
$wsInstalled.Range("C34").Value = "❌"

# --- "Mac installed" sheet: updated notes (column D) ---
$wsInstalled.Range("D17").Value = "zus. Macfuse"
$wsInstalled.Range("D18").Value = "OpenMTP ist besser und wird gepflegt"
$wsInstalled.Range("D24").Value = "05.04.2024: brew installation currently broken; installed manually"

# --- Refresh saved selection/scroll state (matches the last cell the
#     author had selected on each sheet when saving) ---
$wsAvailability.Range("C3").Select() | Out-Null
$wsInstalled.Range("C21").Select() | Out-Null
